$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1 (+3): "- Google Maps API" -> split into "- Google Maps " / "API"
# with the _GoBack bookmark re-inserted between the two runs. Because a
# document can only have one bookmark named "_GoBack", re-adding it here
# implicitly relocates it away from its old spot after "Logout: Youssef"
# (which covers the removal described separately in the diff).
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("- Google Maps API", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $splitPos = $rng1.Start + ("- Google Maps ").Length
    $d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))
}

# ---------------------------------------------------------------------------
# Change 2: "- afterevent forum/private messaging" (with spell-check markup
# around "afterevent") -> "- after" / " " / "event forum/private messaging"
# as three clean runs with no proofErr wrapper.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("- afterevent forum/private messaging", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Text = ""

    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>- after</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>event forum/private messaging</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $ins2 = $d.Range($rng2.Start, $rng2.Start)
    $ins2.InsertXML($xmlFrag)
}
